# Apply cell updates per the crypto price refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.899.15"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "'1.639.39"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'212.51"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "'23.46"
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D10").Value = "'0.0613"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("E11").Value = "  -2.15%  "
$ws.Range("D12").Value = "'1.871.11"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "'1.639.68"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").Value = "'0.563"
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("D17").Value = "'27.883.12"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").Value = "'231.86"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.0₃0723"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'7.65"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "'10.74"
$ws.Range("E22").Value = "  +8.48%  "
$ws.Range("D23").Value = "'4.38"
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("D24").Value = "'2.14"
$ws.Range("D25").Value = "'150.87"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("D26").Value = "'6.91"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").Value = "'0.111"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").Value = "'15.70"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").Value = "'1.18"
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("D33").Value = "'1.458.20"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").Value = "'0.888"
$ws.Range("E37").Value = "  +2.59%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.563"
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0168"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.927"
$ws.Range("E40").Value = "  -2.80%  "
$ws.Range("D41").Value = "'69.19"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'1.02"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").Value = "'2.22"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'1.80"
$ws.Range("E46").Value = "  +6.78%  "
$ws.Range("D47").Value = "'5.35"
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("D48").Value = "'1.781.96"
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("D49").Value = "'88.09"
$ws.Range("E49").Value = "  +2.23%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.100"
$ws.Range("E50").Value = "  +1.83%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.72"
$ws.Range("E51").Value = "  -0.57%  "

Write-Output "Applied 106 cell updates"
